$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Шаблон"

# Clear out old numeric data that is no longer part of the model
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null
$ws.Range("A2").Value = $null

# New header text cell
$ws.Range("A1").Value = "coped cell"

# Column width for column A (raw stored width ends up at 22 after Excel's
# character-width <-> pixel-width conversion)
$ws.Columns.Item(1).ColumnWidth = 21.17

# Row height for header row
$ws.Rows.Item(1).RowHeight = 36.75

# Header style - bold, italic, underline, size 20, Arial Narrow, white font on black fill, centered
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Font.Italic = $true
$header.Font.Underline = $true
$header.Font.Size = 20
$header.Font.Name = "Arial Narrow"
$header.Font.ThemeColor = 2
$header.Interior.ThemeColor = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# Left border on A1 only
$ws.Cells.Item(1, 1).Borders.Item(7).LineStyle = 1
$ws.Cells.Item(1, 1).Borders.Item(7).Weight = 2

# Merge the header cells
$header.Merge()

# Selection matches the merged header
$header.Select()

# Page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
